$d = $word.ActiveDocument

$pairs = @(
    ,@("54-15=", "87-72=")
    ,@("22+16=", "17+36=")
    ,@("42-2=", "11+31=")
    ,@("22+21=", "46+39=")
    ,@("5+91=", "67+16=")
    ,@("32+67=", "64-36=")
    ,@("10+62=", "77-5=")
    ,@("3+4=", "56+24=")
    ,@("85-3=", "91-7=")
    ,@("20+62=", "8-4=")
    ,@("69-29=", "81-30=")
    ,@("84-59=", "50-47=")
    ,@("66-58=", "23-19=")
    ,@("77+21=", "57+11=")
    ,@("56+36=", "43-14=")
    ,@("95-69=", "39+11=")
    ,@("83-50=", "75+1=")
    ,@("91-11=", "90-87=")
    ,@("0+33=", "40-6=")
    ,@("87-23=", "35+29=")
    ,@("49-34=", "77-33=")
    ,@("57+23=", "33-28=")
    ,@("79-55=", "67-41=")
    ,@("47+37=", "64-32=")
    ,@("92-19=", "3+13=")
    ,@("19-11=", "97-33=")
    ,@("26+39=", "23+10=")
    ,@("39-26=", "38+28=")
    ,@("87-43=", "54-47=")
    ,@("91-69=", "83-21=")
    ,@("86-54=", "26-4=")
    ,@("38-27=", "42-34=")
    ,@("15+58=", "76-31=")
    ,@("6+79=", "84-55=")
    ,@("55+12=", "30+12=")
    ,@("49+20=", "51-45=")
    ,@("53-44=", "94-10=")
    ,@("68-8=", "17+1=")
    ,@("42-17=", "9+67=")
    ,@("24-6=", "18+46=")
    ,@("73-3=", "16+48=")
    ,@("40+55=", "87-86=")
    ,@("27+17=", "87-32=")
    ,@("20+6=", "70-13=")
    ,@("10+47=", "95-41=")
    ,@("47-3=", "3+79=")
    ,@("69-59=", "94-53=")
    ,@("41-19=", "40+24=")
    ,@("58+33=", "95-55=")
    ,@("85-8=", "17-7=")
    ,@("69-36=", "79-48=")
    ,@("85-79=", "21+25=")
    ,@("19+40=", "96-36=")
    ,@("67-56=", "67-7=")
    ,@("83-72=", "71-60=")
    ,@("15-3=", "27-14=")
    ,@("29-10=", "19+4=")
    ,@("26-8=", "37-28=")
    ,@("69+14=", "32+20=")
    ,@("24+24=", "44+27=")
    ,@("75-48=", "11+12=")
    ,@("28+71=", "45+3=")
    ,@("89-49=", "95-18=")
    ,@("60+26=", "94+3=")
    ,@("59+7=", "40-10=")
    ,@("68+15=", "39+57=")
    ,@("32+44=", "44-32=")
    ,@("60+3=", "62-7=")
    ,@("7+76=", "75-18=")
    ,@("59+38=", "28+19=")
    ,@("44+52=", "0+4=")
    ,@("65+33=", "60+1=")
    ,@("23-7=", "76-24=")
    ,@("56-0=", "27+9=")
    ,@("30+48=", "31-22=")
    ,@("51-13=", "1+32=")
    ,@("28+63=", "60+31=")
    ,@("65-65=", "36-28=")
    ,@("89-45=", "2+7=")
    ,@("94-35=", "56-16=")
    ,@("34-18=", "57+10=")
    ,@("7+74=", "93-91=")
    ,@("54+20=", "8-7=")
    ,@("42-10=", "44-0=")
    ,@("27+45=", "68-6=")
    ,@("66+17=", "0+55=")
    ,@("75+4=", "2+36=")
    ,@("43-0=", "50-26=")
    ,@("21+28=", "94-42=")
    ,@("33+7=", "32+28=")
    ,@("79-37=", "15-14=")
    ,@("10+8=", "26+63=")
    ,@("66-48=", "62+34=")
    ,@("77-69=", "97-14=")
    ,@("0+20=", "53-1=")
    ,@("83-81=", "86-81=")
    ,@("12+62=", "44-34=")
    ,@("90-6=", "26+64=")
    ,@("52+47=", "78-74=")
    ,@("38+7=", "46+33=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
